$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1/AE1/AF1 need to match the bold/centered/bordered
# style used by the rest of row 1 (style index 1). Copy format from an
# existing header cell rather than re-building font/border/alignment by
# hand so we re-use the existing style instead of creating a new one.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record columns: every data row (2-58) gets the same
# Wins / Losses / Ties values.
for ($r = 2; $r -le 58; $r++) {
    $ws.Cells.Item($r, 30).Value = 76
    $ws.Cells.Item($r, 31).Value = 86
    $ws.Cells.Item($r, 32).Value = 0
}
